$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-report row was inserted at row 473 (pushing the existing
# rows 473..579 down to 474..580). Insert the row first so everything
# below shifts down, then fill in the new row's data.
$ws.Rows.Item(473).Insert()

$ws.Cells.Item(473, 1).Value  = 5
$ws.Cells.Item(473, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(473, 3).Value  = "Maule"
$ws.Cells.Item(473, 4).Value  = 44798
$ws.Cells.Item(473, 5).Value  = 7
$ws.Cells.Item(473, 6).Value  = 100114001
$ws.Cells.Item(473, 7).Value  = "Papa"
$ws.Cells.Item(473, 8).Value  = "Rosara"
$ws.Cells.Item(473, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(473, 10).Value = 1500
$ws.Cells.Item(473, 11).Value = 6000
$ws.Cells.Item(473, 12).Value = 6000
$ws.Cells.Item(473, 13).Value = 6000
$ws.Cells.Item(473, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(473, 15).Value = "Región del Maule"
$ws.Cells.Item(473, 16).Value = 240
$ws.Cells.Item(473, 17).Value = 25
$ws.Cells.Item(473, 18).Value = "Hortaliza"
